# Fix for configurations that have wrong deployed date
# The "Date Deployed" banner in cell A1 was showing the old date
# (18/2/2019); update it to the corrected deployment date (21/2/2019).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Date Deployed: 21/2/2019"

# Reset the view back to the top-left / default selection.
$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
